$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 472. Excel shifts rows 472:587 down to 473:588,
# and the new blank row 472 inherits formatting (incl. the date style on column D)
# from the row it was inserted above.
$ws.Rows(472).Insert()

# The columns that stay identical to the row immediately below (which now holds
# what used to be row 472's data) are copied across so the new record matches
# the same market/product/variety/quality combination.
$ws.Range("A472").Value = $ws.Range("A473").Value()
$ws.Range("B472").Value = $ws.Range("B473").Value()
$ws.Range("C472").Value = $ws.Range("C473").Value()
$ws.Range("E472").Value = $ws.Range("E473").Value()
$ws.Range("F472").Value = $ws.Range("F473").Value()
$ws.Range("G472").Value = $ws.Range("G473").Value()
$ws.Range("H472").Value = $ws.Range("H473").Value()
$ws.Range("I472").Value = $ws.Range("I473").Value()
$ws.Range("J472").Value = $ws.Range("J473").Value()
$ws.Range("K472").Value = $ws.Range("K473").Value()
$ws.Range("L472").Value = $ws.Range("L473").Value()

# New record's own data (week of 2023-10-12 at Región de O'Higgins).
$ws.Range("D472").Value = 45211
$ws.Range("M472").Value = 180
$ws.Range("N472").Value = 15000
$ws.Range("O472").Value = 15000
$ws.Range("P472").Value = 15000
$ws.Range("Q472").Value = "$/bandeja 18 kilos granel"
$ws.Range("R472").Value = "Región de O'Higgins"
$ws.Range("S472").Value = 833
$ws.Range("T472").Value = 18
